$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("Objetivos:") - B/C now hold the responsible student's name
$ws.Range("B10:C10").Value = "7290967 - Emerson Gonçalves de Melo"

# Row 13 ("Programa resumido:") - B/C now hold the activation date text.
# Assigning the literal "01/01/2023" string directly would be auto-parsed
# by Excel into a date serial number (and force a new number format),
# so we write it as a formula returning the text and then flatten the
# formula down to a static value via copy / paste-special-values. This
# keeps the cell as a shared-string text cell using the existing style.
$ws.Range("B13:C13").Formula = '="01/01/2023"'
$ws.Range("B13:C13").Copy()
$ws.Range("B13:C13").PasteSpecial(-4163)

# Row 15 ("Programa:") - B/C now hold the same responsible student's name
$ws.Range("B15:C15").Value = "7290967 - Emerson Gonçalves de Melo"

# Row 18 ("Método:") - B/C now hold the other responsible professor's name
$ws.Range("B18:C18").Value = "1176388 - Luiz Tadeu Fernandes Eleno"
